$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Annotations")

# Insert 3 new rows at row 6, pushing existing rows 6-14 down to 9-17.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(8).Insert()

# Fill in the new rows with the added fields.
$ws.Cells.Item(6,1).Value = "~descr"
$ws.Cells.Item(6,2).Value = "The Description field of the Event"
$ws.Cells.Item(6,3).Value = "Gas 10"

$ws.Cells.Item(7,1).Value = "~type"
$ws.Cells.Item(7,2).Value = "The Type field of the Event"
$ws.Cells.Item(7,3).Value = "Power"

$ws.Cells.Item(8,1).Value = "~sldrunit"
$ws.Cells.Item(8,2).Value = "The value of the Slider Unit for this Event"
$ws.Cells.Item(8,3).Value = "kPa"

# Sheet view updates: Annotations becomes the selected/active tab, scrolled to top,
# with A9 selected.
$ws.Activate() | Out-Null
$ws.Range("A9").Select() | Out-Null
